$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "319.90"
Set-TextValue "E2" "3.55%"
Set-TextValue "G2" "20"
Set-TextValue "D3" "41.40"
Set-TextValue "E3" "1.12%"
Set-TextValue "G3" "20"
Set-TextValue "D4" "5.257"
Set-TextValue "E4" "2.59%"
Set-TextValue "G4" "20"
Set-TextValue "D5" "0.07745"
Set-TextValue "E5" "1.66%"
Set-TextValue "G5" "20"
Set-TextValue "D6" "1.741"
Set-TextValue "E6" "8.46%"
Set-TextValue "G6" "20"
Set-TextValue "D7" "0.9459"
Set-TextValue "E7" "4.17%"
Set-TextValue "G7" "20"
Set-TextValue "E8" "-2.36%"
Set-TextValue "G8" "20"
Set-TextValue "D9" "0.1240"
Set-TextValue "E9" "-1.48%"
Set-TextValue "G9" "20"
Set-TextValue "D10" "0.1866"
Set-TextValue "E10" "3.37%"
Set-TextValue "G10" "20"
Set-TextValue "D11" "0.09149"
Set-TextValue "E11" "1.02%"
Set-TextValue "G11" "20"
Set-TextValue "D12" "0.04136"
Set-TextValue "E12" "-3.80%"
Set-TextValue "G12" "20"
Set-TextValue "E13" "0.59%"
Set-TextValue "G13" "20"
Set-TextValue "D14" "0.001288"
Set-TextValue "E14" "2.81%"
Set-TextValue "G14" "20"
Set-TextValue "D15" "0.005843"
Set-TextValue "E15" "3.28%"
Set-TextValue "G15" "20"
Set-TextValue "G16" "20"
Set-TextValue "D17" "3.355"
Set-TextValue "E17" "0.09%"
Set-TextValue "G17" "20"
Set-TextValue "D18" "4.347"
Set-TextValue "E18" "1.50%"
Set-TextValue "G18" "20"
Set-TextValue "E19" "1.36%"
Set-TextValue "G19" "20"
Set-TextValue "D20" "8.755"
Set-TextValue "E20" "26.56%"
Set-TextValue "G20" "20"
Set-TextValue "D21" "0.1353"
Set-TextValue "E21" "-2.89%"
Set-TextValue "G21" "20"
Set-TextValue "D22" "0.2824"
Set-TextValue "E22" "3.04%"
Set-TextValue "G22" "20"
Set-TextValue "D23" "0.04037"
Set-TextValue "E23" "0.01%"
Set-TextValue "G23" "20"
Set-TextValue "E24" "0.09%"
Set-TextValue "G24" "20"
Set-TextValue "D25" "0.004114"
Set-TextValue "E25" "1.91%"
Set-TextValue "G25" "20"
Set-TextValue "E26" "-0.23%"
Set-TextValue "G26" "20"
Set-TextValue "G27" "20"
Set-TextValue "G28" "20"
Set-TextValue "G29" "20"
Set-TextValue "G30" "20"
Set-TextValue "G31" "20"
Set-TextValue "G32" "20"
Set-TextValue "G33" "20"
Set-TextValue "G34" "20"
Set-TextValue "G35" "20"
Set-TextValue "G36" "20"
Set-TextValue "G37" "20"
Set-TextValue "D38" "0.02562"
Set-TextValue "E38" "6.14%"
Set-TextValue "G38" "20"
Set-TextValue "D39" "0.05336"
Set-TextValue "E39" "2.04%"
Set-TextValue "G39" "20"
Set-TextValue "D40" "0.007753"
Set-TextValue "E40" "-1.17%"
Set-TextValue "G40" "20"
Set-TextValue "D41" "0.1317"
Set-TextValue "E41" "1.14%"
Set-TextValue "G41" "20"
Set-TextValue "D42" "0.007046"
Set-TextValue "E42" "3.61%"
Set-TextValue "G42" "20"
Set-TextValue "D43" "0.001992"
Set-TextValue "E43" "8.09%"
Set-TextValue "G43" "20"
Set-TextValue "D44" "0.008250"
Set-TextValue "E44" "10.51%"
Set-TextValue "G44" "20"
Set-TextValue "D45" "0.3177"
Set-TextValue "E45" "-5.51%"
Set-TextValue "G45" "20"
Set-TextValue "D46" "0.00006697"
Set-TextValue "E46" "-2.69%"
Set-TextValue "G46" "20"
Set-TextValue "E47" "-0.22%"
Set-TextValue "G47" "20"
Set-TextValue "D48" "0.2009"
Set-TextValue "E48" "49.39%"
Set-TextValue "G48" "20"
Set-TextValue "D49" "0.004204"
Set-TextValue "E49" "40.00%"
Set-TextValue "G49" "20"
Set-TextValue "E50" "-0.22%"
Set-TextValue "G50" "20"
Set-TextValue "E51" "-0.22%"
Set-TextValue "G51" "20"
